# Enabled Equalizer in codec
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Date:" header to the new date (2018-02-01, serial 43132)
$ws.Range("C1").Value = "2/1/2018"

# Row 34: "Status LEDs" task is now Closed
$ws.Range("C34").Value = "Closed"

# Row 62: "Root cause voice_data message loss and corruption in duplex mode" is now Closed
$ws.Range("C62").Value = "Closed"

# New row 69: follow-up task for voice data retransmission
$ws.Range("B69").Value = "Implement voice data retransmission protocol"
$ws.Range("C69").Value = "Open"

# Row 55: rename task and mark as Closed (equalizer feature shipped)
$ws.Range("B55").Value = "Use 5 band equalizer"
$ws.Range("C55").Value = "Closed"

# Row 56: "Install stronger speaker" moved from Open to Ongoing
$ws.Range("C56").Value = "Ongoing"

# Row 57: "Only send EchoReq as keep-alive when no incoming comm." moved from Open to Ongoing
$ws.Range("C57").Value = "Ongoing"

# New row 70: newly logged bug
$ws.Range("B70").Value = "Bug: sometimes msg comm goes haywire after speech"
$ws.Range("C70").Value = "Open"

# Update selection / active cell and scroll position to match the new view state
$ws.Activate()
try {
    $excel.ActiveWindow.ScrollRow = 37
    $excel.ActiveWindow.ScrollColumn = 1
} catch { }
$ws.Range("C58").Select()
